# Re-orders entries inside several python-list-style text cells
# (columns F/G/H = "Matched Udsteder" / "Matched Værdipapirets navn" / "Kommuner")
# for a handful of rows, matching updated upstream matching/export logic.
# The set of list elements in each touched cell is unchanged - only the
# ordering within the cell's text differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - Carnival Corp
$ws.Range("F17").Value = '[''Carnival Corp'', nan, ''Carnival Corp.'', ''USA'', ''CARNIVAL CORP'', ''Carnival Corp 15.01.2030'']'
$ws.Range("G17").Value = '[''Carnival Corp'', ''Carnival Corp.'', ''9,875 Carnival Corp 01-08-2027'', ''7.63% Carnival Corp 3/2026'', ''CARNIVAL CORP 5.75% 01.03.2027'', ''CARNIVAL CORP 4% 01.08.2028'', ''6 CCL 01-05-2029 (144A)'', ''CARNIVAL CORP 6% 01.05.2029'', ''CARNIVAL CORP 7% 15.08.2029'', ''5,75 CCL 15-01-2030 (REGS)'', ''5.75% Carnival Corp 2030'', ''CARNIVAL CORP 5.75% 15.01.2030'', ''CCL 5 3/4 01/15/30'']'
$ws.Range("H17").Value = '[''Ikast-Brande'', ''Faaborg-Midtfyn'', ''Guldborgsund'', ''Hedensted'', ''Helsingør'', ''Hørsholm'', ''Middelfart'', ''Norddjurs'', ''Nordfyn'', ''Roskilde'', ''Slagelse'', ''Svendborg'', ''Syddjurs'', ''Esbjerg'', ''Odense'', ''Rødovre'', ''Fanø'', ''Aarhus'', ''Herning'', ''Vejen'', ''Kalundborg'', ''Ballerup'', ''Region Hovedstaden'', ''Brøndby'', ''Brønderslev'', ''Hvidovre'', ''Kerteminde'', ''Kolding'', ''Køge'', ''Lolland'', ''Region Midtjylland'', ''Region Sjælland'', ''Silkeborg'', ''Skive'', ''Sorø'', ''Varde'', ''Fredensborg'', ''Struer'', ''Region Nordjylland'', ''Fredericia'', ''Greve'', ''Høje Taastrup'', ''Lemvig'', ''Næstved'', ''Randers'', ''Hillerød'', ''Aabenraa'']'
# Row 82 - Target Corp
$ws.Range("F82").Value = '[''Target\xa0Corp'', nan, ''Target Corp'', ''Target Corp.'', ''USA'']'
$ws.Range("G82").Value = '[''Target'', ''TARGET CORP'', ''Target Corp'', ''Target Corp.'']'
$ws.Range("H82").Value = '[''Viborg'', ''Assens'', ''Mariagerfjord'', ''Ballerup'', ''Brøndby'', ''Brønderslev'', ''Fredensborg'', ''Frederiksberg'', ''Faaborg-Midtfyn'', ''Hedensted'', ''Herning'', ''Holstebro'', ''Hvidovre'', ''Ikast-Brande'', ''Kerteminde'', ''Kolding'', ''Køge'', ''Lolland'', ''Morsø'', ''Region Midtjylland'', ''Region Nordjylland'', ''Region Sjælland'', ''Rudersdal'', ''Silkeborg'', ''Skive'', ''Sorø'', ''Struer'', ''Varde'', ''Odense'', ''Aalborg'', ''Guldborgsund'', ''Helsingør'', ''Hørsholm'', ''Middelfart'', ''Norddjurs'', ''Nordfyn'', ''Roskilde'', ''Slagelse'', ''Svendborg'', ''Syddjurs'', ''Rødovre'', ''Fanø'', ''Esbjerg'', ''Aarhus'']'
# Row 87 - Vale SA
$ws.Range("G87").Value = '[''VALE SA'', ''Vale SA'', ''VALE SA-SP ADR'']'
# Row 91 - Walmart Inc
$ws.Range("F91").Value = '[''WALMART INC'', ''Walmart\xa0Inc'', ''Wal-Mart Stores'', ''Walmart Inc.'']'
$ws.Range("G91").Value = '[''Walmart'', ''WALMART INC'', ''Walmart Inc'', ''Walmart Inc.'', ''Wal-Mart Stores'', ''Wal-Mart Stores Inc'', ''Walmart, Inc.'']'
$ws.Range("H91").Value = '[''Guldborgsund'', ''Næstved'', ''Odense'', ''Randers'', ''Region\xa0Nordjylland'', ''Viborg'', ''Mariagerfjord'', ''Skive'', ''Svendborg'', ''Fredericia'', ''Greve'', ''Herning'', ''Hillerød'', ''Høje Taastrup'', ''Hørsholm'', ''Kolding'', ''Lemvig'', ''Varde'', ''Aabenraa'', ''Egedal'', ''Favrskov'', ''Region Sjælland'', ''Fanø'', ''Assens'', ''Billund'', ''Bornholm'', ''Frederiksberg'', ''Gribskov'', ''Holstebro'', ''Horsens'', ''Hvidovre'', ''Jammerbugt'', ''Lejre'', ''Lyngby-Taarbæk'', ''Nordfyn'', ''Nyborg'', ''Rebild'', ''Roskilde'', ''Rødovre'', ''Skanderborg'', ''Solrød'', ''Thisted'', ''Tønder'', ''Vallensbæk'', ''Norddjurs'', ''Aalborg'', ''Faaborg-Midtfyn'', ''Hedensted'', ''Helsingør'', ''Middelfart'', ''Slagelse'', ''Syddjurs'', ''Esbjerg'', ''Aarhus'']'
